$d = $word.ActiveDocument

# Commit message: "change time table poster session starts at 9:00"
# The poster-session row currently reads "10:00-12:00 (Tehran)" /
# "5:30-7:30 (UTC)". It should become "9:00-12:00 (Tehran)" /
# "4:30-7:30 (UTC)".

# Move/remove the existing "_GoBack" tracking bookmark first (Word
# relocates this automatically to the point of the most recent edit),
# then relocate it to the new edit point after the text is updated.
$d.Bookmarks.ShowHidden = $true
$hadGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hadGoBack) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$d.Content.Find.Execute("10:00-12:00 (Tehran)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "9:00-12:00 (Tehran)", 2)

$d.Content.Find.Execute("5:30-7:30 (UTC)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4:30-7:30 (UTC)", 2)

# Re-anchor "_GoBack" right after the "4" that replaced the old "5", which
# is where Word leaves it once the last keystroke of the edit landed.
$target = $d.Content
$target.Find.Execute("4:30-7:30 (UTC)")
$splitPos = $target.Start + 1
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)
